$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 283
    $ws.Range("F3").Value = 226
    $ws.Range("F4").Value = 27

    # G4 holds its number-looking value as text (matches the other cells in
    # column G, e.g. G2/G3/G5), so force text entry, write it, then drop the
    # leftover explicit cell style back to Normal so only the value/type
    # changes (no visual formatting change remains on the cell).
    $ws.Range("G4").NumberFormat = "@"
    $ws.Range("G4").Value = "45"
    $ws.Range("G4").Style = "Normal"

    $ws.Range("F5").Value = 265
}
